# This script recreates the edit described by the diff:
# - Adds three new weekly blocks (rows 19-27, 29-37, 39-47) after the existing
#   two weekly blocks, plus a fourth block (rows 48-56) continuing the pattern,
#   each block following the same "title / header / 4 member rows / blank / summary"
#   layout as the existing blocks (rows 1-9 and 10-18).
# - Re-uses the existing visual styles (bold title bars, bordered header row,
#   bordered body rows, left-aligned summary bars) by copying formatting from
#   the already existing analogous rows.
# - Introduces one new font/style for the "plan content" column text in the
#   new blocks, and one new font/style for the new bold, centered date titles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Application.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Prepare the two new template styles used by the new blocks.
# ---------------------------------------------------------------------------

# New bold / centered "date title" style (e.g. A19, A29, A39, A48).
# Start from the existing plain title style (A1) and tweak it.
Copy-Format "A1:D1" "A19"
$ws.Range("A19").Font.Name = "宋体"
$ws.Range("A19").Font.Bold = $true
$ws.Range("A19").Font.Size = 10
$ws.Range("A19").HorizontalAlignment = $xlCenter

# New "plan content" column style used in column B of the new blocks.
# Start from the existing plain body style and tweak the font so it becomes
# a distinct style entry. Use a single source cell so the paste does not
# tile across multiple columns.
Copy-Format "B3" "B21"
$ws.Range("B21").Font.Name = "宋体"
$ws.Range("B21").Font.Size = 11
$ws.Range("B21").Font.Bold = $false
$ws.Range("B21").Value = ""

# ---------------------------------------------------------------------------
# Helper that builds one complete weekly block starting at row $r:
#   r+0 : merged bold date title                (A)
#   r+1 : header row 组员/计划内容/完成情况/备注 (A:D)
#   r+2..r+5 : four member rows                  (A:D)
#   r+6 : blank separator row                    (A:D)
#   r+7,r+8 : merged "总结：" bar                 (A:D)
# ---------------------------------------------------------------------------
function Build-Block($r, $dateText, $members, $tasks, $statuses) {
    $title = $r
    $header = $r + 1
    $row1 = $r + 2
    $row2 = $r + 3
    $row3 = $r + 4
    $row4 = $r + 5
    $blank = $r + 6
    $sumTop = $r + 7
    $sumBot = $r + 8

    # Title bar: merge first, then copy format uniformly, then set text.
    $ws.Range("A$($title):D$($title)").Merge() | Out-Null
    Copy-Format "A19:D19" "A$($title):D$($title)"
    $ws.Range("A$title").Value = $dateText

    # Header row.
    Copy-Format "A2:D2" "A$($header):D$($header)"
    $ws.Range("A$header").Value = "组员"
    $ws.Range("B$header").Value = "计划内容"
    $ws.Range("C$header").Value = "完成情况"
    $ws.Range("D$header").Value = "备注"

    # Four member rows.
    $dataRows = @($row1, $row2, $row3, $row4)
    for ($i = 0; $i -lt 4; $i++) {
        $rr = $dataRows[$i]
        Copy-Format "A3:D3" "A$($rr):D$($rr)"
        Copy-Format "B21" "B$($rr)"
        $ws.Range("A$rr").Value = $members[$i]
        $ws.Range("B$rr").Value = $tasks[$i]
        $ws.Range("C$rr").Value = $statuses[$i]
        $ws.Range("D$rr").Value = ""
    }

    # Blank separator row (still bordered).
    Copy-Format "A7:D7" "A$($blank):D$($blank)"

    # Summary bar.
    $ws.Range("A$($sumTop):D$($sumBot)").Merge() | Out-Null
    Copy-Format "A8:D9" "A$($sumTop):D$($sumBot)"
    $ws.Range("A$sumTop").Value = "总结："
}

# ---------------------------------------------------------------------------
# Block 1: rows 19-27 - 日期：2018.10.10 第六周周三
# ---------------------------------------------------------------------------
Build-Block 19 "日期：2018.10.10 第六周周三" `
    @("余舒章","王嘉宇","许俊杰","庞森杰") `
    @("编写分配到的用例规约","编写分配到的用例规约","编写分配到的用例规约","编写分配到的用例规约") `
    @("已完成","已完成","已完成","已完成")

# ---------------------------------------------------------------------------
# Block 2: rows 29-37 - 日期：2018.10.11 第六周周四
# ---------------------------------------------------------------------------
Build-Block 29 "日期：2018.10.11 第六周周四" `
    @("余舒章","王嘉宇","许俊杰","庞森杰") `
    @("修改及完善用例规约","修改及完善用例规约","修改及完善用例规约","设计安卓端界面") `
    @("已完成","已完成","已完成","进行中")

# ---------------------------------------------------------------------------
# Block 3: rows 39-47 - 日期：2018.10.15 第七周周一
# ---------------------------------------------------------------------------
Build-Block 39 "日期：2018.10.15 第七周周一" `
    @("余舒章","王嘉宇","许俊杰","庞森杰") `
    @("讨论并设计er图","讨论并设计er图","讨论并设计er图","设计安卓端界面") `
    @("进行中","进行中","进行中","进行中")

# ---------------------------------------------------------------------------
# Block 4: rows 48-56 - 日期：2018.10.17 第七周周三
# ---------------------------------------------------------------------------
Build-Block 48 "日期：2018.10.17 第七周周三" `
    @("余舒章","王嘉宇","许俊杰","庞森杰") `
    @("讨论并设计er图","讨论并设计er图","讨论并设计er图","设计安卓端界面") `
    @("进行中","进行中","进行中","进行中")

# ---------------------------------------------------------------------------
# Final view state: scrolled down to the newly added content.
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A19"), $true) | Out-Null
$ws.Range("D51").Select() | Out-Null
